$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Capture row 16 and row 17 current values (B:G) before swapping
$b16 = $ws.Range("B16").Value2
$c16 = $ws.Range("C16").Value2
$d16 = $ws.Range("D16").Value2
$e16 = $ws.Range("E16").Value2
$f16 = $ws.Range("F16").Value2
$g16 = $ws.Range("G16").Value2

$b17 = $ws.Range("B17").Value2
$c17 = $ws.Range("C17").Value2
$d17 = $ws.Range("D17").Value2
$e17 = $ws.Range("E17").Value2
$f17 = $ws.Range("F17").Value2
$g17 = $ws.Range("G17").Value2

# Write row17's original data into row16, and row16's original data into row17 (swap)
$ws.Range("B16").Value = $b17
$ws.Range("C16").Value = $c17
$ws.Range("D16").Value = $d17
$ws.Range("E16").Value = $e17
$ws.Range("F16").Value = $f17
$ws.Range("G16").Value = $g17

$ws.Range("B17").Value = $b16
$ws.Range("C17").Value = $c16
$ws.Range("D17").Value = $d16
$ws.Range("E17").Value = $e16
$ws.Range("F17").Value = $f16
$ws.Range("G17").Value = $g16
